# Update the Quiz slide's question text (slide 6 of the deck).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$rsquo = [char]8217

$tr.Paragraphs(1,1).Runs(1,1).Text = "Subclass objects can" + $rsquo + "t override some of the template method call.  (True/False)"
$tr.Paragraphs(2,1).Runs(1,1).Text = "The template method cannot be define in the abstract class.           (True/False)"
$tr.Paragraphs(3,1).Runs(1,1).Text = "Why don" + $rsquo + "t we need template method if the subclasses will not be making similar calls?"
